$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info: cardholder name and card number ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number stored as TEXT in the original file.
# Force text storage (apostrophe prefix) then restore the original
# (non quote-prefixed) formatting by copying format from a plain text
# neighbor cell, so the resulting style matches the untouched style id.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance label ---
$ws.Range("D5").Value = "KONTOSTAND AM 17.05.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "19.05."
$ws.Range("C6").Value = "20.05."
$ws.Range("D6").Value = "PAYPAL LATQYQ"
$ws.Range("E6").Value = "60,91-"

# --- Row 7 ---
$ws.Range("B7").Value = "20.05."
$ws.Range("C7").Value = "21.05."
$ws.Range("D7").Value = "BURGER KING Sebnitz"
$ws.Range("E7").Value = "44,99-"

# --- Row 8 ---
$ws.Range("B8").Value = "23.05."
$ws.Range("C8").Value = "24.05."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-18609504"
$ws.Range("E8").Value = "56,16-"

# --- Row 9 ---
$ws.Range("B9").Value = "26.05."
$ws.Range("C9").Value = "27.05."
$ws.Range("D9").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E9").Value = "59,51-"

# --- Row 10 (previously an empty placeholder row, now filled in) ---
$ws.Range("B10").Value = "30.05."
$ws.Range("C10").Value = "31.05."
$ws.Range("D10").Value = "KARTENZ./30.05 LIDL RO"
$ws.Range("E10").Value = "56,64-"
# E10 switches from the wrapped/centered placeholder style to the same
# right-aligned style used by the other amount cells (E6:E9).
$ws.Range("E10").HorizontalAlignment = $ws.Range("E9").HorizontalAlignment
$ws.Range("E10").VerticalAlignment = $ws.Range("E9").VerticalAlignment
$ws.Range("E10").WrapText = $ws.Range("E9").WrapText

# --- Row 11 (previously an empty placeholder row, now filled in) ---
$ws.Range("B11").Value = "31.05."
$ws.Range("C11").Value = "01.06."
$ws.Range("D11").Value = "KARTENZ./31.05 ALDI SUED RO"
$ws.Range("E11").Value = "47,64-"
$ws.Range("E11").HorizontalAlignment = $ws.Range("E9").HorizontalAlignment
$ws.Range("E11").VerticalAlignment = $ws.Range("E9").VerticalAlignment
$ws.Range("E11").WrapText = $ws.Range("E9").WrapText

# --- Closing balance ---
$ws.Range("D12").Value = "KONTOSTAND AM 05.06.2024"
$ws.Range("E12").Value = "325,85-"

# --- Next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.06.2024"
